$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (StreetTRACKS Gold Shares / GLD) - was Newmont Corporation / NEM
$ws.Range("B2").Value = "StreetTRACKS Gold Shares"
$ws.Range("C2").Value = "GLD"
$ws.Range("D2").Value = 387.88
$ws.Range("F2").Value = 3.48
$ws.Range("H2").Value = 56
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 68.8
$ws.Range("N2").Value = 85.96878041621773

# Row 3 (Newmont Corporation / NEM) - was StreetTRACKS Gold Shares / GLD
$ws.Range("B3").Value = "Newmont Corporation"
$ws.Range("C3").Value = "NEM"
$ws.Range("D3").Value = 91
$ws.Range("F3").Value = 11.29
$ws.Range("H3").Value = 66
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 68.8
$ws.Range("N3").Value = 85.96878041621773

# Row 4 (Gold Dec 25 / GC=F) - updated score values
$ws.Range("K4").Value = 63
$ws.Range("N4").Value = 85.96878041621773
